$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet (end of tab strip).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ShopForChargers"

# Reuse the existing "description" cell format (Menlo font, s=2) from the
# ShopForDSLRs sheet so no new style/font entries get created.
$fmtSource = $wb.Worksheets.Item("ShopForDSLRs").Range("C2")
$fmtSource.Copy($ws.Range("B2:D2"))

# Values for the new "shop for chargers" search-result row.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "7"
$ws.Range("B2").Value = "Dodge"
$ws.Range("C2").Value = "Charger"
$ws.Range("D2").Value = "Dodge Charger"

# Column D needs to be a bit wider to fit "Dodge Charger".
$ws.Columns.Item(4).ColumnWidth = 16.5

# Make the new sheet the active tab with D2 selected (mirrors the other
# "shop for X" sheets, which keep the last populated cell selected).
$wb.Worksheets.Item("ShopForChargers").Activate()
[void]$ws.Range("D2").Select()

# Window shrank slightly between saves.
$excel.ActiveWindow.Height = 16060
